# Update EC database: replace worker LAUREANO HERRERA MELENDEZ (73140681) on
# row 16 with VANESSA DE LA PEÑA VILLALBA (1143348947), drop the extra
# CRISTIAN ANDRES LIMA PUERTA period row, and renumber/ascend the period
# column for VANESSA from 2002 to 2104 (part 1 of the new estado de cuenta).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove two of the duplicate "VANESSA" period rows (rows 30 and 31 both
# already belong to VANESSA) - this shifts the CRISTIAN row (32) up to 30,
# and the signature block (37/38) up to (35/36), matching the new layout.
$ws.Rows("30:31").Delete()

# --- Header figures -----------------------------------------------------
$ws.Range("E11").Value = 4200000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 15

# --- Detail rows 16-29: same worker (VANESSA), ascending period ---------
$periods = @("2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("C$r").Value = "1143348947"
    $ws.Range("D$r").Value = "VANESSA DE LA PEÑA VILLALBA"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = 280000
    $ws.Range("G$r").Value = 7000000
}

# --- Final detail row (was the CRISTIAN ANDRES LIMA PUERTA row, now at 30)
$ws.Range("C30").Value = "1143348947"
$ws.Range("D30").Value = "VANESSA DE LA PEÑA VILLALBA"
$ws.Range("E30").Value = "2104"
$ws.Range("F30").Value = 280000
$ws.Range("G30").Value = 7000000

# --- Cosmetic: column D auto-fit changed slightly with the new text -----
$ws.Columns("D").ColumnWidth = 29.26953125
